$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.791.53'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '2.087.68'
$ws.Range('E3').Value = '  +4.67%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.66'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.58%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.387'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0767'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('D13').Value = '2.392.84'
$ws.Range('E13').Value = '  +4.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.785'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.54%  '
$ws.Range('D18').Value = '2.097.70'
$ws.Range('E18').Value = '  +5.02%  '
$ws.Range('D19').Value = '37.829.23'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +20.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').Value = '0.0₃0817'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.133'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('E31').Value = '  +6.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0633'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.33%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.03%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.25%  '
$ws.Range('B39').Value = 'BinanceUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0969'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.54%  '
$ws.Range('D43').Value = '1.487.86'
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.43%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0213'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.15%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.35%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +26.73%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.61%  '
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.06%  '
